$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.935.72"
$ws.Range("E2").Value = "  +0.28%  "

$ws.Range("D3").Value = "1.646.58"
$ws.Range("E3").Value = "  +1.03%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.28%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.40"
$ws.Range("E5").Value = "  +0.34%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5063"
$ws.Range("E6").Value = "  -0.11%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.005"
$ws.Range("E7").Value = "  +0.31%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2592"
$ws.Range("E8").Value = "  +0.66%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06451"
$ws.Range("E9").Value = "  +1.93%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.64"
$ws.Range("E10").Value = "  +6.09%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07824"
$ws.Range("E11").Value = "  +0.87%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.288"
$ws.Range("E12").Value = "  +1.02%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.644.57"
$ws.Range("E13").Value = "  +0.72%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "1.869.59"
$ws.Range("E14").Value = "  +0.80%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5628"
$ws.Range("E15").Value = "  +2.37%  "

$ws.Range("D16").Value = "0.0₅7732"
$ws.Range("E16").Value = "  +0.98%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.49"
$ws.Range("E17").Value = "  -0.32%  "

$ws.Range("D18").Value = "25.958.51"
$ws.Range("E18").Value = "  +0.32%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.005"
$ws.Range("E19").Value = "  +0.32%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "194.20"
$ws.Range("E20").Value = "  -0.12%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.381"
$ws.Range("E21").Value = "  -0.55%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.968"
$ws.Range("E22").Value = "  +0.98%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.142"
$ws.Range("E23").Value = "  +1.92%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.005"
$ws.Range("E24").Value = "  +0.24%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.802"
$ws.Range("E25").Value = "  -5.88%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.47"
$ws.Range("E26").Value = "  -0.43%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1236"
$ws.Range("E27").Value = "  -0.98%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.832"
$ws.Range("E28").Value = "  +1.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.58"
$ws.Range("E29").Value = "  +0.08%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.247"
$ws.Range("E30").Value = "  +0.74%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04972"
$ws.Range("E31").Value = "  +1.73%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.305"
$ws.Range("E32").Value = "  +1.73%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.246"
$ws.Range("E33").Value = "  +1.77%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.581"
$ws.Range("E34").Value = "  +2.46%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.389"
$ws.Range("E35").Value = "  +0.83%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9082"
$ws.Range("E36").Value = "  +1.46%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5577"
$ws.Range("E37").Value = "  +1.02%  "

$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.577"
$ws.Range("E38").Value = "  +1.59%  "

$ws.Range("D39").Value = "1.133.55"
$ws.Range("E39").Value = "  +1.35%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01568"
$ws.Range("E40").Value = "  +1.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.005"
$ws.Range("E41").Value = "  +0.43%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.526"
$ws.Range("E42").Value = "  -1.01%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8033"
$ws.Range("E43").Value = "  +0.81%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "98.78"
$ws.Range("E44").Value = "  +1.60%  "

$ws.Range("D45").Value = "1.778.41"
$ws.Range("E45").Value = "  +0.71%  "

$ws.Range("D46").Value = "0.0₈109"
$ws.Range("E46").Value = "  -8.57%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.85"
$ws.Range("E47").Value = "  +2.10%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4287"
$ws.Range("E48").Value = "  -3.55%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.771"
$ws.Range("E49").Value = "  +2.76%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05048"
$ws.Range("E50").Value = "  -1.69%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9986"
$ws.Range("E51").Value = "  -0.69%  "
